$wb = $excel.ActiveWorkbook

# --- Sheet "Overview" (sheet1) -------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E1:F2").EntireColumn.ColumnWidth = 29.9777047293527

# --- Sheet "zh-cn" (sheet2) ------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K2").Value = "2016-08-31 10:54:55"
$wsZhCn.Range("P2").Value = ""
$wsZhCn.Range("C1:C2").EntireColumn.ColumnWidth = 29.9777047293527
$wsZhCn.Range("P1:P2").EntireColumn.ColumnWidth = 13.7470528738839

# --- Sheet "de-de" (sheet3) ------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K2").Value = "2016-08-31 10:55:11"
$wsDeDe.Range("P2").Value = ""
$wsDeDe.Range("C1:C2").EntireColumn.ColumnWidth = 29.9777047293527
$wsDeDe.Range("P1:P2").EntireColumn.ColumnWidth = 13.7470528738839
